# Update mods data [2026-01-21 15:19:01]
# Append a new data row (row 72) to the ModCounts sheet with the latest
# mod-count reading for 逃离鸭科夫.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = 72
$sourceRow = 71

# Force column A to be treated as plain text so the date-like string
# "2026/01/21" is not auto-converted into a date serial number (the rest
# of the column stores these as literal text, not real dates).
$ws.Cells.Item($newRow, 1).NumberFormat = "@"

$ws.Cells.Item($newRow, 1).Value = "2026/01/21"
$ws.Cells.Item($newRow, 2).Value = "逃离鸭科夫"
$ws.Cells.Item($newRow, 3).Value = 1153

# Match the formatting (center/center alignment, default number format)
# used by the rest of the data rows by copying the format from the row
# directly above.
$ws.Range("A" + $sourceRow + ":C" + $sourceRow).Copy($null)
$ws.Range("A" + $newRow + ":C" + $newRow).PasteSpecial(-4122)

$excel.CutCopyMode = $false
